$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typos in the rating-description text blocks (column D) ---

# D4: "unable to peform any one elses duties" -> "unable to perform any one else's duties"
$d4 = @'
5: Demonstrates KSAs to do excellent work, acquires new KSA to help team, can perform any role on team if necessary
4: Between 5 above and 3 below
3: Demonstrates sufficient KSA to contribute to team, acquires KSAs to meet requirements, able to perform other tasks
2: Between 3 above and 1 below
1: Missing basic qualification, unable to develop KSAs to contribute to team, unable to perform any one else's duties
'@
$ws.Range("D4").Value = $d4

# D6: "contributsions" -> "contributions"
$d6 = @'
5: Is interested in teammates ideas and contributions, makes sure everyone is informed, is encouraging, enthusiastic and asks for feedback/suggestions
3: Listens and respects teammate contributions, communicates clearly, shares info, participates fully, reacts and responds to feedback/suggestions
1: Interrupts, ignores, bosses, or makes fun, takes action without input, does not share, complains, makes excuses, does not interact, is defensive
'@
$ws.Range("D6").Value = $d6

# D7: "sucess" -> "success"
$d7 = @'
5: Monitors teams' progress, makes sure teammates are progressing, gives specific, timely, and constructive feedback
3: Knows what everyone on the team should be doing and notices problems, alerts teammates and suggests solutions with success is threatened
1: Unaware if team is meeting goals, does not pay attention to teammates progress, avoids discussing team problems even when obvious
'@
$ws.Range("D7").Value = $d7

# D8: "Encouarges" -> "Encourages", "responsiblities" -> "responsibilities"
$d8 = @'
5: Motivates team to do excellent work, cares about excellent work even without reward, believes in team's ability to do excellent work
3: Encourages good work to meet requirements, believes team can meet its responsibilities
1: Satisfied even if not all requirements are met,  avoids work, doubts team can meet requirements
'@
$ws.Range("D8").Value = $d8

# D3: "for each each member" -> "for each member"
$ws.Range("D3").Value = "Rating Descriptions (provide whole number ratings (5, 4, 3, 2, or 1) in columns for each member including yourself)"

# --- Clean up unused/blank styled cells trailing in rows 9-11 ---
$ws.Range("B9:C9").Clear()
$ws.Range("E9:U9").Clear()
$ws.Range("F10:U10").Clear()
$ws.Range("A11").Clear()
$ws.Range("F11:U11").Clear()

# E10/E11 lose their explicit (redundant) cell style, reverting to the default style
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Style = "Normal"

# --- Update the active selection to D3 (was A4) ---
$ws.Range("D3").Select()
